# Applies the "cryptos list" data refresh described in the commit:
#   "Updated cryptos list on Sat Sep 16 20:24:27 UTC 2023 with GitHub Actions"
#
# Mostly Price (column D) / Volume(1h) (column E) updates, plus rows 49/50
# which swap position: BabyDogeCoin <-> EnergySwap (name, link, price, volume).
#
# Note: several new Price values are plain numeric-looking strings (e.g. "216.24").
# Excel auto-converts those to numbers (losing the original text formatting /
# trailing zeros), so we prefix those specific assignments with a leading "'"
# (PowerShell single-quote escape '''...') to force Excel to keep them as text,
# matching the source data which stores every column as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '26.738.59'
$ws.Range("E2").Value = '  +0.84%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.646.28'
$ws.Range("E3").Value = '  +1.18%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.41%  '

# Row 5: BNB
$ws.Range("D5").Value = '''216.24'
$ws.Range("E5").Value = '  +1.45%  '

# Row 6: XRP
$ws.Range("E6").Value = '  +0.41%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.34%  '

# Row 8: Cardano
$ws.Range("E8").Value = '  +0.88%  '

# Row 9: Dogecoin
$ws.Range("E9").Value = '  +0.45%  '

# Row 10: Solana
$ws.Range("E10").Value = '  +2.21%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.0843'
$ws.Range("E11").Value = '  -0.21%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.875.93'
$ws.Range("E12").Value = '  +1.20%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.669.38'
$ws.Range("E13").Value = '  +2.28%  '

# Row 14: Polkadot
$ws.Range("E14").Value = '  +1.13%  '

# Row 15: Polygon
$ws.Range("E15").Value = '  +1.66%  '

# Row 16: Litecoin
$ws.Range("D16").Value = '''65.33'
$ws.Range("E16").Value = '  +0.32%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '26.747.28'
$ws.Range("E17").Value = '  +0.76%  '

# Row 18: ShibaInu
$ws.Range("E18").Value = '  +0.36%  '

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''218.24'
$ws.Range("E19").Value = '  +1.67%  '

# Row 20: Dai
$ws.Range("E20").Value = '  +0.35%  '

# Row 21: Uniswap
$ws.Range("E21").Value = '  +1.78%  '

# Row 22: Toncoin
$ws.Range("E22").Value = '  +16.20%  '

# Row 23: Chainlink
$ws.Range("D23").Value = '''6.26'
$ws.Range("E23").Value = '  -0.04%  '

# Row 24: Avalanche
$ws.Range("D24").Value = '''9.48'
$ws.Range("E24").Value = '  +1.88%  '

# Row 25: Monero
$ws.Range("D25").Value = '''146.51'
$ws.Range("E25").Value = '  -0.79%  '

# Row 26: BinanceUSD
$ws.Range("E26").Value = '  +0.38%  '

# Row 27: Stellar
$ws.Range("E27").Value = '  -0.06%  '

# Row 28: Cosmos
$ws.Range("D28").Value = '''7.16'
$ws.Range("E28").Value = '  +3.99%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''15.73'
$ws.Range("E29").Value = '  +1.23%  '

# Row 30: Hedera
$ws.Range("E30").Value = '  +1.59%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  +1.55%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  -0.05%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = '  +1.23%  '

# Row 34: Maker
$ws.Range("D34").Value = '1.279.12'
$ws.Range("E34").Value = '  +2.92%  '

# Row 35: LidoDAOToken
$ws.Range("E35").Value = '  +2.99%  '

# Row 36: HuobiToken
$ws.Range("E36").Value = '  +3.02%  '

# Row 37: VeChain
$ws.Range("E37").Value = '  +2.13%  '

# Row 38: ImmutableX
$ws.Range("D38").Value = '''0.537'
$ws.Range("E38").Value = '  +5.52%  '

# Row 39: ARBITRUM
$ws.Range("E39").Value = '  +4.38%  '

# Row 40: PaxDollar
$ws.Range("E40").Value = '  +0.34%  '

# Row 41: TrustWalletToken
$ws.Range("E41").Value = '  +2.13%  '

# Row 42: MXToken
$ws.Range("D42").Value = '''2.24'
$ws.Range("E42").Value = '  -1.09%  '

# Row 43: FraxShare
$ws.Range("D43").Value = '''5.44'
$ws.Range("E43").Value = '  +1.97%  '

# Row 44: RocketPoolETH
$ws.Range("D44").Value = '1.786.78'
$ws.Range("E44").Value = '  +1.32%  '

# Row 45: Quant
$ws.Range("D45").Value = '''92.04'
$ws.Range("E45").Value = '  -1.29%  '

# Row 46: Aave
$ws.Range("D46").Value = '''59.70'
$ws.Range("E46").Value = '  +8.81%  '

# Row 47: RenderToken
$ws.Range("E47").Value = '  +1.46%  '

# Row 48: Cronos
$ws.Range("E48").Value = '  +1.11%  '

# Row 49: was BabyDogeCoin, now EnergySwap
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.75'
$ws.Range("E49").Value = '  +3.28%  '

# Row 50: was EnergySwap, now BabyDogeCoin
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₇0985'
$ws.Range("E50").Value = '  -5.30%  '

# Row 51: Algorand
$ws.Range("E51").Value = '  +1.42%  '
